# Generate Report for Handoff
# Adds a new row (row 3) to each of the three worksheets (Overview, zh-cn,
# de-de) reflecting a new file "d19248fe-1236-4c5c-8b9d-ba696d9793d3o...md"
# that has reached "Ready for handoff" status, resizes the tables to include
# the new row, and adds hyperlinks pointing at the new source file.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

$mdName      = "d19248fe-1236-4c5c-8b9d-ba696d9793d3ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$mdNameE2E   = "e2e\d19248fe-1236-4c5c-8b9d-ba696d9793d3ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$xlfZhCn     = "d19248fe-1236-4c5c-8b9d-ba696d9793d3oooooooooooooooooooooooooooooooooooooooo.19eb863f006c003738113f568e75324a6232a1e6.zh-cn.xlf"
$xlfDeDe     = "d19248fe-1236-4c5c-8b9d-ba696d9793d3oooooooooooooooooooooooooooooooooooooooo.19eb863f006c003738113f568e75324a6232a1e6.de-de.xlf"
$readyStatus = "Ready for handoff"
$dtHandoffZh = "2016-08-24 00:27:15"
$dtHandoffDe = "2016-08-24 00:27:20"
$dtNull      = "0001-01-01 00:00:00"

$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c9582e2317521d07586830e0f9864d962f61390e/e2e/" + $mdName

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Range("B3").Value = $mdNameE2E
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $dtHandoffDe
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $mdNameE2E) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)

$wsZhCn.Range("A3").Value = $mdName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = $dtHandoffZh
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("K3").Value = $dtNull
$wsZhCn.Range("K3").NumberFormat = $dateFmt
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkUrl, "", "", $mdName) | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)

$wsDeDe.Range("A3").Value = $mdName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = $dtHandoffDe
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("K3").Value = $dtNull
$wsDeDe.Range("K3").NumberFormat = $dateFmt
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkUrl, "", "", $mdName) | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

# ---------------------------------------------------------------------
# Column width tweaks (status columns widen to fit "Ready for handoff")
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3

Write-Output "Report generated for handoff"
